# edit.ps1
# Applies the "Updated cryptos list" refresh to the crypto price sheet.
# For each changed row, updates the Price (column D) and Volume(1h) (column E)
# cells with their new text values, preserving the original inline-string /
# text cell typing (no numeric reinterpretation, no stray style changes).

function Set-TextValue {
    param($cell, [string]$text)
    # If the text would be auto-parsed by Excel as a plain number (e.g. "243.11"),
    # force it to stay text by using the quote-prefix, then strip the resulting
    # quotePrefix style so the cell keeps the default (unstyled) formatting.
    if ($text -match '^[+-]?\d+(\.\d+)?$') {
        $cell.Value = "'" + $text
        $cell.Style = "Normal"
    } else {
        $cell.Value = $text
    }
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

Set-TextValue $ws.Range("D2") '42.237.85'
Set-TextValue $ws.Range("E2") '  -0.62%  '
Set-TextValue $ws.Range("D3") '2.235.16'
Set-TextValue $ws.Range("E3") '  -0.15%  '
Set-TextValue $ws.Range("E4") '  -0.07%  '
Set-TextValue $ws.Range("D5") '243.11'
Set-TextValue $ws.Range("E5") '  -0.62%  '
Set-TextValue $ws.Range("D6") '0.628'
Set-TextValue $ws.Range("E6") '  -0.08%  '
Set-TextValue $ws.Range("D7") '74.68'
Set-TextValue $ws.Range("E7") '  -0.34%  '
Set-TextValue $ws.Range("E8") '  +0.04%  '
Set-TextValue $ws.Range("E9") '  -2.38%  '
Set-TextValue $ws.Range("D10") '42.58'
Set-TextValue $ws.Range("E10") '  -2.17%  '
Set-TextValue $ws.Range("E11") '  +1.65%  '
Set-TextValue $ws.Range("E12") '  -2.55%  '
Set-TextValue $ws.Range("E13") '  +1.03%  '
Set-TextValue $ws.Range("D14") '2.569.33'
Set-TextValue $ws.Range("E14") '  -0.26%  '
Set-TextValue $ws.Range("D15") '14.38'
Set-TextValue $ws.Range("E15") '  -0.79%  '
Set-TextValue $ws.Range("E16") '  -1.79%  '
Set-TextValue $ws.Range("D17") '2.249.27'
Set-TextValue $ws.Range("E17") '  -0.85%  '
Set-TextValue $ws.Range("D18") '42.110.16'
Set-TextValue $ws.Range("E18") '  -0.41%  '
Set-TextValue $ws.Range("E19") '  +2.51%  '
Set-TextValue $ws.Range("E20") '  +1.01%  '
Set-TextValue $ws.Range("D21") '72.90'
Set-TextValue $ws.Range("E21") '  +1.38%  '
Set-TextValue $ws.Range("D22") '11.38'
Set-TextValue $ws.Range("E22") '  +9.24%  '
Set-TextValue $ws.Range("D23") '230.79'
Set-TextValue $ws.Range("E23") '  -0.09%  '
Set-TextValue $ws.Range("E24") '  -5.61%  '
Set-TextValue $ws.Range("E25") '  -0.04%  '
Set-TextValue $ws.Range("D26") '11.38'
Set-TextValue $ws.Range("E26") '  -2.26%  '
Set-TextValue $ws.Range("E27") '  -0.14%  '
Set-TextValue $ws.Range("D28") '2.28'
Set-TextValue $ws.Range("E28") '  -0.93%  '
Set-TextValue $ws.Range("D29") '2.20'
Set-TextValue $ws.Range("E29") '  -2.75%  '
Set-TextValue $ws.Range("D30") '167.44'
Set-TextValue $ws.Range("E30") '  +0.39%  '
Set-TextValue $ws.Range("D31") '20.64'
Set-TextValue $ws.Range("E31") '  -0.13%  '
Set-TextValue $ws.Range("E32") '  -4.25%  '
Set-TextValue $ws.Range("E33") '  -1.16%  '
Set-TextValue $ws.Range("D34") '30.05'
Set-TextValue $ws.Range("E34") '  +0.18%  '
Set-TextValue $ws.Range("E35") '  -0.28%  '
Set-TextValue $ws.Range("D36") '0.110'
Set-TextValue $ws.Range("E36") '  -6.85%  '
Set-TextValue $ws.Range("E37") '  -5.17%  '
Set-TextValue $ws.Range("E38") '  -2.62%  '
Set-TextValue $ws.Range("D39") '13.28'
Set-TextValue $ws.Range("E39") '  -1.26%  '
Set-TextValue $ws.Range("E40") '  -1.39%  '
Set-TextValue $ws.Range("D41") '5.73'
Set-TextValue $ws.Range("E41") '  +0.68%  '
Set-TextValue $ws.Range("D42") '65.16'
Set-TextValue $ws.Range("E42") '  +2.92%  '
Set-TextValue $ws.Range("E43") '  -0.18%  '
Set-TextValue $ws.Range("E44") '  -0.99%  '
Set-TextValue $ws.Range("D45") '104.95'
Set-TextValue $ws.Range("E45") '  -1.28%  '
Set-TextValue $ws.Range("E46") '  -2.13%  '
Set-TextValue $ws.Range("E47") '  -0.68%  '
Set-TextValue $ws.Range("E48") '  -1.83%  '
Set-TextValue $ws.Range("E49") '  -0.20%  '
Set-TextValue $ws.Range("E50") '  -1.51%  '
Set-TextValue $ws.Range("D51") '2.441.15'
Set-TextValue $ws.Range("E51") '  -0.64%  '
